# Fixed a bug with ONEExtractor at end of month, added debug_mode configuration
#
# On the "Hamburg" configuration sheet, insert a new "debug_mode" (boolean,
# default FALSE) setting right after "mask_date_if_bol_present", and relocate
# "g2_whitespace_rows" so that it immediately follows "randomiser_upper_interval"
# (instead of trailing at the very end), pushing "true_sample"/"false_sample"
# down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hamburg")

# Capture the current (pre-edit) values of the rows that are shifting, before
# any of them get overwritten.
$randLowerLabel = $ws.Cells.Item(13, 1).Value2
$randLowerValue = $ws.Cells.Item(13, 2).Value2
$randUpperLabel = $ws.Cells.Item(14, 1).Value2
$randUpperValue = $ws.Cells.Item(14, 2).Value2
$trueSampleLabel = $ws.Cells.Item(15, 1).Value2
$trueSampleValue = $ws.Cells.Item(15, 2).Value2
$falseSampleLabel = $ws.Cells.Item(16, 1).Value2
$falseSampleValue = $ws.Cells.Item(16, 2).Value2
$whitespaceLabel = $ws.Cells.Item(17, 1).Value2
$whitespaceValue = $ws.Cells.Item(17, 2).Value2

# Rewrite rows 13-18 bottom-up so nothing is clobbered before it's read.

# Row 18: false_sample (was row 16) - brand-new row, needs formats applied.
$ws.Cells.Item(18, 1).Value = $falseSampleLabel
$ws.Cells.Item(18, 2).Value = [bool]$falseSampleValue

# Row 17: true_sample (was row 15)
$ws.Cells.Item(17, 1).Value = $trueSampleLabel
$ws.Cells.Item(17, 2).Value = [bool]$trueSampleValue

# Row 16: g2_whitespace_rows (was row 17) - numeric value, moved up.
$ws.Cells.Item(16, 1).Value = $whitespaceLabel
$ws.Cells.Item(16, 2).Value = $whitespaceValue

# Row 15: randomiser_upper_interval (was row 14) - unchanged position content-wise
$ws.Cells.Item(15, 1).Value = $randUpperLabel
$ws.Cells.Item(15, 2).Value = $randUpperValue

# Row 14: randomiser_lower_interval (was row 13)
$ws.Cells.Item(14, 1).Value = $randLowerLabel
$ws.Cells.Item(14, 2).Value = $randLowerValue

# Row 13: debug_mode (brand new)
$ws.Cells.Item(13, 1).Value = "debug_mode"
$ws.Cells.Item(13, 2).Value = $false

# --- Formatting cleanup -----------------------------------------------
# Column A keeps the same "label" formatting throughout (style already in
# place for every row that previously held data); only the brand-new row 18
# needs it copied in explicitly.
$ws.Cells.Item(12, 1).Copy() | Out-Null
$ws.Cells.Item(18, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Column B alternates between a "number" look (left aligned) and a "boolean"
# look (center aligned). Re-apply the correct one to every touched cell by
# copying formats from a cell that already has the right look.
$ws.Cells.Item(12, 2).Copy() | Out-Null                 # boolean-styled cell
$ws.Cells.Item(13, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(18, 2).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(14, 2).Copy() | Out-Null                 # number-styled cell
$ws.Cells.Item(15, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(16, 2).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Restore the active selection to match the post-edit workbook state.
$ws.Range("B12").Select()
